# Update 'threads' up to 64
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Threaded_Pi_Calc")

# Row 1: thread counts, extend B1:G1 -> B1:I1 with 32, 64
$ws.Range("H1").Value = 32
$ws.Range("I1").Value = 64

# Row 2: win32api timings, updated values extended to H2:I2
$ws.Range("B2").Value = 0.397928
$ws.Range("C2").Value = 0.200374
$ws.Range("D2").Value = 0.098641
$ws.Range("E2").Value = 0.062549
$ws.Range("F2").Value = 0.053941
$ws.Range("G2").Value = 0.051897
$ws.Range("H2").Value = 0.048013
$ws.Range("I2").Value = 0.052322

# Row 3: openmp timings, updated values extended to H3:I3
$ws.Range("B3").Value = 0.39557
$ws.Range("C3").Value = 0.200206
$ws.Range("D3").Value = 0.105281
$ws.Range("E3").Value = 0.056742
$ws.Range("F3").Value = 0.042005
$ws.Range("G3").Value = 0.044003
$ws.Range("H3").Value = 0.04202
$ws.Range("I3").Value = 0.042011

# New columns H:K get the same custom width as the rest of the data columns
$ws.Range("H1:K1").ColumnWidth = 14.140625

# Update sheet view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 205
$ws.Range("I6").Select()
